$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Missing1")

# --- Missing1: selection / active-tab changes (tabSelected moves to the new sheet) ---
$ws1.Range("C7:F13").Select()

# --- Add the new "Missing2" worksheet right after "Missing1" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Missing2"

# Row 5
$ws2.Range("C5").Value = "A"
$ws2.Range("D5").Value = "B"
$ws2.Range("E5").Value = "C"
$ws2.Range("F5").Value = "D"

# Row 6
$ws2.Range("C6").Value = "NA"
$ws2.Range("D6").Value = "a"
$ws2.Range("E6").Formula = "=TRUE"
$ws2.Range("F6").Value = 29921
$ws2.Range("F6").NumberFormat = "mm-dd-yy"

# Row 7
$ws2.Range("C7").Value = -3.2
$ws2.Range("D7").Value = "empty"
$ws2.Range("E7").Value = -9999
$ws2.Range("F7").Value = 29922
$ws2.Range("F7").NumberFormat = "mm-dd-yy"

# Row 8
$ws2.Range("C8").Value = -9999
$ws2.Range("D8").Value = "c"
$ws2.Range("E8").Formula = "=FALSE"
$ws2.Range("F8").Value = "empty"
$ws2.Range("F8").NumberFormat = "mm-dd-yy"

# Row 9 (has the "missing value" style tweak: general format explicitly re-applied)
$ws2.Range("C9").Value = "missing"
$ws2.Range("D9").Value = "x"
$ws2.Range("E9").Value = "missing"
$ws2.Range("F9").Value = -9999

# Row 10
$ws2.Range("C10").Value = 8
$ws2.Range("D10").Value = "a"
$ws2.Range("E10").Formula = "=FALSE"
$ws2.Range("F10").Value = "NA"
$ws2.Range("F10").NumberFormat = "mm-dd-yy"

# Row 11
$ws2.Range("C11").Value = "empty"
$ws2.Range("D11").Value = "o"
$ws2.Range("E11").Value = "NA"
$ws2.Range("F11").Value = 29926
$ws2.Range("F11").NumberFormat = "mm-dd-yy"

$ws2.Range("F9").Select()
